# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# This script:
#   1. Clears the placeholder/empty cells in "ODI Batting Extra" (columns B-F)
#      that never actually held scraped data, so they drop out of the sheet
#      (mirrors cells that simply were never written by the scraper).
#   2. Adds a new "ODI Bowling Extra" worksheet (sheetId 5) after
#      "ODI Batting Extra", with the scraped MAIDEN_OVERS /
#      PERCENT_WICKETS_OF_ALL columns keyed by MATCH_CODE.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Trim the empty placeholder cells out of "ODI Batting Extra"
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$emptyRanges = @(
    "B2:E2", "B3:E3", "B5:E5", "B6:E6", "B7:E7", "E8", "B9:E9",
    "C11:E11", "B12:E12", "B13:E13",
    "B14:F14", "B15:F15", "B16:F16", "B17:F17", "B18:F18", "B19:F19", "B20:F20", "B21:F21"
)
foreach ($rng in $emptyRanges) {
    $battingExtra.Range($rng).Value = ""
}

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Bowling Extra" worksheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row - match the look of the other "*Extra" sheets (bold, boxed,
# centered/top-aligned header cells).
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le 3; $col++) {
    $bowlingExtra.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$headerRange = $bowlingExtra.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data rows: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
# (all three columns are text-typed, so numeric-looking values like "4430"
# or "2" and percentages like "10.00%" are forced to text the same way the
# scraper originally wrote them).
$data = @(
    @(2,  "4430", "2",  "10.00%"),
    @(3,  "4431", "0",  ""),
    @(4,  "4435", "",   ""),
    @(5,  "4436", "",   ""),
    @(6,  "4437", "1",  "10.00%"),
    @(7,  "4483", "1",  "30.00%"),
    @(8,  "4486", "",   ""),
    @(9,  "4594", "",   ""),
    @(10, "4597", "1",  ""),
    @(11, "4600", "",   ""),
    @(12, "4601", "0",  ""),
    @(13, "4603", "3",  "20.00%"),
    @(14, "4644", "1",  ""),
    @(15, "4645", "",   ""),
    @(16, "4646", "2",  "30.00%"),
    @(17, "4647", "",   ""),
    @(18, "4648", "2",  ""),
    @(19, "4649", "0",  "10.00%"),
    @(20, "4663", "",   ""),
    @(21, "4666", "",   "")
)

foreach ($row in $data) {
    $r = $row[0]
    $matchCode = $row[1]
    $maidenOvers = $row[2]
    $percentWickets = $row[3]

    $codeCell = $bowlingExtra.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $matchCode
    $codeCell.ClearFormats()

    $oversCell = $bowlingExtra.Cells.Item($r, 2)
    $oversCell.NumberFormat = "@"
    $oversCell.Value = $maidenOvers
    $oversCell.ClearFormats()

    $pctCell = $bowlingExtra.Cells.Item($r, 3)
    $pctCell.NumberFormat = "@"
    $pctCell.Value = $percentWickets
    $pctCell.ClearFormats()
}
